$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear row 1 entirely (contents + formatting)
$ws.Rows("1:1").Clear()

# Clear columns D:G entirely (contents + formatting)
$ws.Range("D:G").Clear()

# Type the new value into B8
$ws.Range("B8").Value = "Punam Bhoyar"

# Set the selection to B8, matching the final state
$ws.Range("B8").Select()
